# Atualização de bases das ligas, do dia: 19-06-2024 às 21:51
#
# The underlying match records (columns B:AD) for several rows in the
# "Peru Liga 1" sheet were re-ordered (the sequential "id" in column A
# stays tied to the row, but the actual match data that belongs to it
# moved to a different row). This script snapshots the B:AD data for
# every affected row and then writes it back out according to the
# required new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps: destination row number -> source row number (the row whose
# original B:AD content must end up at the destination row).
$rowMap = @{
    156 = 157
    157 = 156
    184 = 185
    185 = 186
    186 = 184
    187 = 188
    188 = 187
    305 = 306
    306 = 305
    312 = 313
    313 = 312
    338 = 340
    339 = 338
    340 = 339
}

# Snapshot the current B:AD values of every row referenced above before
# writing anything back out, so that rows which both give and receive
# data (e.g. the 156/157 swap, or the 184/185/186 three-way rotation)
# are handled correctly regardless of write order.
$snapshot = @{}
foreach ($srcRow in ($rowMap.Values | Sort-Object -Unique)) {
    $snapshot[$srcRow] = $ws.Range("B$srcRow`:AD$srcRow").Value2
}

foreach ($destRow in ($rowMap.Keys | Sort-Object)) {
    $srcRow = $rowMap[$destRow]
    $ws.Range("B$destRow`:AD$destRow").Value2 = $snapshot[$srcRow]
}
